$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "EmployeeCode" (B) and "Email" (C) columns
$ws.Range("B1:C1").EntireColumn.Delete()

# Add new "RoleId" header and values
$ws.Range("D1").Value = "RoleId"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3

$ws.Range("F4").Select()
